# Update production_date values in column D (rows 2-6) as per commit:
# "Update test data to use relative dates instead of fixed future dates"
# The cells store the dates as plain text strings (inline strings), so we
# force a text number format first to avoid Excel auto-converting the
# assigned string into a date serial value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "2026-02-12"
    "D3" = "2026-02-13"
    "D4" = "2026-02-14"
    "D5" = "2026-02-15"
    "D6" = "2026-02-16"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
